# Apply the commit's changes:
#  - sheet1 is rewritten with new row data (rows 1-3), row 11 untouched
#  - "sheet2" is duplicated from sheet1 (keeps styles/namespaces) and placed
#    right after it, then given its own row data
#  - sheet2 becomes the active sheet/tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Duplicate sheet1 -> sheet2, inserted immediately after sheet1 ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "sheet2"

# --- Rewrite sheet1 data ---
$ws1.Range("A1").Value = 44
$ws1.Range("B1").Value = 45
$ws1.Range("C1").Value = "third"

$ws1.Range("A2").Value = 34
$ws1.Range("B2").Value = 11.0
$ws1.Range("C2").Value = 12.0

$ws1.Range("A3").Value = 21
$ws1.Range("B3").Value = 21.0
$ws1.Range("C3").Value = 22.0

# row 11 (B11=99) stays as-is

# Selection on sheet1 becomes the A1:C3 block
[void]$ws1.Range("A1:C3").Select()

# --- Rewrite sheet2 data ---
$ws2.Range("B11").ClearContents()

$ws2.Range("A1").Value = 31.0
$ws2.Range("B1").Value = 32.0
$ws2.Range("C1").Value = "third"

$ws2.Range("A2").Value = 41.0
$ws2.Range("B2").Value = 42.0
$ws2.Range("C2").Value = 12

$ws2.Range("A3").Value = 21
$ws2.Range("B3").Value = 21
$ws2.Range("C3").Value = 22

# sheet2 is the active sheet/tab, with B8 selected
$ws2.Activate()
[void]$ws2.Range("B8").Select()
